$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Apply the "Table Grid" table style (adds <w:tblStyle w:val="TableGrid"/> to tblPr)
$t.Style = "Table Grid"

# Update the "t" (column 12) and "t_max" (column 13) values for each data row.
# Row 2 .. Row 7 correspond to std = 0.0, 0.1, 0.2, 0.3, 0.5, 1.0
$t.Cell(2, 12).Range.Text = "0.0027"
$t.Cell(2, 13).Range.Text = "0.0053"

$t.Cell(3, 12).Range.Text = "0.0028"
$t.Cell(3, 13).Range.Text = "0.0055"

$t.Cell(4, 12).Range.Text = "0.0028"
$t.Cell(4, 13).Range.Text = "0.0068"

$t.Cell(5, 12).Range.Text = "0.0028"
$t.Cell(5, 13).Range.Text = "0.0057"

$t.Cell(6, 12).Range.Text = "0.0028"
$t.Cell(6, 13).Range.Text = "0.0059"

$t.Cell(7, 12).Range.Text = "0.0028"
$t.Cell(7, 13).Range.Text = "0.0054"

# Append a new, completely empty row at the bottom of the table.
$newRow = $t.Rows.Add()
$newRowIndex = $newRow.Index
for ($col = 1; $col -le $t.Columns.Count; $col++) {
    $t.Cell($newRowIndex, $col).Range.Delete()
}
